$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.022.28"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "3.408.77"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'575.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'126.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.60%  "
$ws.Range("D8").Value = "3.409.95"
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("D9").Value = "'0.476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("E12").Value = "  -1.79%  "
$ws.Range("D13").Value = "3.998.11"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("E14").Value = "  -0.80%  "
$ws.Range("D15").Value = "3.418.56"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("E16").Value = "  -4.21%  "
$ws.Range("D17").Value = "62.991.62"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").Value = "'24.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.15%  "
$ws.Range("D19").Value = "'9.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("D20").Value = "'5.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").Value = "'13.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").Value = "'377.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.04%  "
$ws.Range("D23").Value = "'0.559"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("D24").Value = "3.550.58"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").Value = "'72.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.59%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("E27").Value = "  -7.48%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -5.37%  "
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("D31").Value = "'7.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.61%  "
$ws.Range("E32").Value = "  -3.27%  "
$ws.Range("E33").Value = "  -4.79%  "
$ws.Range("D34").Value = "3.443.33"
$ws.Range("E34").Value = "  -2.48%  "
$ws.Range("D36").Value = "'22.75"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("D37").Value = "'5.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").Value = "'6.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("D39").Value = "'164.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("E40").Value = "  -3.75%  "
$ws.Range("D41").Value = "'0.0762"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'0.780"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").Value = "'41.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("E46").Value = "  -5.03%  "
$ws.Range("D47").Value = "'22.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.48%  "
$ws.Range("E48").Value = "  -7.23%  "
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").Value = "2.260.02"
$ws.Range("E50").Value = "  -5.27%  "
$ws.Range("E51").Value = "  -3.96%  "
